$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the "Name" (column A) values.
#     Edit the Skill* rows first, then the NormalAtk*/NormalThump rows, so
#     new shared-string entries are appended in the same order the author's
#     session produced them (SKILL1..4 before NORMALATTACK1..3/NORMALTHUMP).
$ws.Range("A6").Value = "SKILL1"
$ws.Range("A7").Value = "SKILL2"
$ws.Range("A8").Value = "SKILL3"
$ws.Range("A9").Value = "SKILL4"

$ws.Range("A2").Value = "NORMALATTACK1"
$ws.Range("A3").Value = "NORMALATTACK2"
$ws.Range("A4").Value = "NORMALATTACK3"
$ws.Range("A5").Value = "NORMALTHUMP"

# --- Every row's "NextLevelID" (column C) now points at NORMALATTACK2 ---
$ws.Range("C2").Value = "NORMALATTACK2"
$ws.Range("C3").Value = "NORMALATTACK2"
$ws.Range("C4").Value = "NORMALATTACK2"
$ws.Range("C5").Value = "NORMALATTACK2"
$ws.Range("C6").Value = "NORMALATTACK2"
$ws.Range("C7").Value = "NORMALATTACK2"
$ws.Range("C8").Value = "NORMALATTACK2"
$ws.Range("C9").Value = "NORMALATTACK2"

# --- Bump the AnimaState id column (G) for the Skill rows ---
$ws.Range("G6").Value = 101
$ws.Range("G7").Value = 102
$ws.Range("G8").Value = 103
$ws.Range("G9").Value = 104

# --- Format the whole NextLevelID column's data as Text, matching the
#     plain (no banding fill/border) look already used by the table rows ---
$ws.Range("C2:C9").NumberFormat = "@"
$ws.Range("C2").Copy()
$ws.Range("C4:C9").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Reflect the author's final selection (NextLevelID column data) ---
$ws.Range("C2:C9").Select() | Out-Null
